$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 18) for "Brass Birmingham Board Game"
$ws.Range("A18").Value = "Brass Birmingham Board Game "
$ws.Range("B18").Value = "https://m.media-amazon.com/images/I/614zbkJJJgL._AC_SL1200_.jpg"
$ws.Range("C18").Value = "https://www.amazon.de/-/en/Roxley-Games-Brass-Birmingham-Board/dp/1988884047/ref=sr_1_1?crid=2DVJHJTB01TNV&keywords=Brass%2BBirmingham&qid=1699391429&s=toys&sprefix=brass%2Bbirmingham%2Ctoys%2C100&sr=1-1&th=1"
$ws.Range("D18").Value = "94.31 EUR"

# Move the selection as recorded in the saved workbook
$ws.Range("D19").Select()
